$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used throughout:
#  - Some new text values look like dates/numbers ("01/01/2023", "2") and Excel's
#    smart-entry would silently convert them to a date serial / numeric value. To
#    keep them as plain text (as in the source data) we temporarily force Text
#    number format before assigning, then restore the row/column's normal look by
#    copying the format from an always-text reference cell (B2/C2) that keeps the
#    same style classes (s="2" for column B, s="3" for column C) used everywhere
#    else on this sheet.
#  - Some target cells are brand new (didn't exist in the original row) and this
#    runtime gives them the wrong default style when first created, so the same
#    format-copy step is applied there as well.

# ---- Row 5: Créditos-aula: 4 -> 2 ----
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null

# ---- Row 7: Carga horária: 60 h -> 30 h (stays text naturally) ----
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# ---- Row 8: Ativação: 01/01/2012 -> 01/01/2023 ----
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "01/01/2023"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null

# ---- Row 10: Objetivos responsible professor (stays text naturally) ----
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# ---- Row 11: Objectives - new English objectives text (brand new cells) ----
$ws.Range("B11").Value = "To present notions of fluid mechanics, through the study of fluid media when static or in motion. Enable the student to model and solve problems of interest in fluid mechanics, with adequate choice of hypotheses and application of corresponding solution tools."
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null

$ws.Range("C11").Value = "To present notions of fluid mechanics, through the study of fluid media when static or in motion. Enable the student to model and solve problems of interest in fluid mechanics, with adequate choice of hypotheses and application of corresponding solution tools."
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null

# ---- Row 13: Programa resumido value (mirrors Ativação date, per source data) ----
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "01/01/2023"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

# ---- Row 14: Short syllabus - new English text (brand new cells) ----
$ws.Range("B14").Value = "Fundamentals of fluid mechanics. Introduction to fluid statics. Integral and differential formulation of mass, energy and momentum transport equations. Dimensional analysis and similarity. Incompressible flow of ideal and viscous fluids, laminar and turbulent regime. Navier-Stokes equation. Boundary layer theory."
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null

$ws.Range("C14").Value = "Fundamentals of fluid mechanics. Introduction to fluid statics. Integral and differential formulation of mass, energy and momentum transport equations. Dimensional analysis and similarity. Incompressible flow of ideal and viscous fluids, laminar and turbulent regime. Navier-Stokes equation. Boundary layer theory."
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null

# ---- Row 15: Programa value (mirrors Objetivos professor, per source data; stays text naturally) ----
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"

# ---- Row 16: Syllabus - new English full text (brand new cells) ----
$ws.Range("B16").Value = "Introduction: fluid concept; properties and concept of continuum; modeling of transfer processes; analysis methods; dimensions and units. • Fluid statics review: basic hydrostatic equation, pressure variation in a static fluid; principles of Stevin, Pascal and Archimedes. • Integral formulation of transport equations: Reynolds transport theorem; application to the principles of conservation of mass, momentum and energy; Bernoulli equation. • Differential formulation of transport equations: description of the flow; differential form: from the principles of conservation of mass, momentum and energy; dimensionless formulation, dimensional analysis and similarity. Dimensionless groups: Reynolds number and Grashoff number. • Internal incompressible flow: Euler equations; Newton's law for viscosity, shear stresses; Navier-Stokes equation; flow regimes: laminar and turbulent flow. Calculation of pressure drop (distributed and localized), friction coefficient. • External incompressible flow: introduction to the boundary layer; flow around bodies, drag force."
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null

$ws.Range("C16").Value = "Introduction: fluid concept; properties and concept of continuum; modeling of transfer processes; analysis methods; dimensions and units. • Fluid statics review: basic hydrostatic equation, pressure variation in a static fluid; principles of Stevin, Pascal and Archimedes. • Integral formulation of transport equations: Reynolds transport theorem; application to the principles of conservation of mass, momentum and energy; Bernoulli equation. • Differential formulation of transport equations: description of the flow; differential form: from the principles of conservation of mass, momentum and energy; dimensionless formulation, dimensional analysis and similarity. Dimensionless groups: Reynolds number and Grashoff number. • Internal incompressible flow: Euler equations; Newton's law for viscosity, shear stresses; Navier-Stokes equation; flow regimes: laminar and turbulent flow. Calculation of pressure drop (distributed and localized), friction coefficient. • External incompressible flow: introduction to the boundary layer; flow around bodies, drag force."
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null

# ---- Row 18: Método - new docente responsible name (stays text naturally) ----
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$excel.CutCopyMode = 0
